$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '35.003.58'
$ws.Range('E2').Value = '  +0.16%  '
$ws.Range('D3').Value = '1.826.94'
$ws.Range('E3').Value = '  +0.54%  '
$ws.Range('E4').Value = '  +0.19%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '230.99'
$ws.Range('E5').Value = '  -0.42%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.617'
$ws.Range('E6').Value = '  +1.00%  '
$ws.Range('E7').Value = '  +0.12%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '40.09'
$ws.Range('E8').Value = '  -2.44%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.322'
$ws.Range('E9').Value = '  +5.56%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0683'
$ws.Range('E10').Value = '  -0.05%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0994'
$ws.Range('E11').Value = '  -0.68%  '
$ws.Range('D12').Value = '2.091.52'
$ws.Range('E12').Value = '  +0.50%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '11.33'
$ws.Range('E13').Value = '  +2.58%  '
$ws.Range('E14').Value = '  +1.69%  '
$ws.Range('D15').Value = '1.820.35'
$ws.Range('E15').Value = '  +0.16%  '
$ws.Range('E16').Value = '  +0.04%  '
$ws.Range('D17').Value = '35.089.96'
$ws.Range('E17').Value = '  +0.48%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '69.72'
$ws.Range('E18').Value = '  +0.68%  '
$ws.Range('E19').Value = '  +0.15%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '240.69'
$ws.Range('E20').Value = '  +1.22%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.11'
$ws.Range('E21').Value = '  +3.49%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.74'
$ws.Range('E22').Value = '  +2.17%  '
$ws.Range('E23').Value = '  -0.01%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.26'
$ws.Range('E24').Value = '  +0.09%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '174.06'
$ws.Range('E25').Value = '  +0.91%  '
$ws.Range('E26').Value = '  +0.66%  '
$ws.Range('E27').Value = '  +3.38%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '17.36'
$ws.Range('E28').Value = '  -0.20%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.52'
$ws.Range('E29').Value = '  -4.84%  '
$ws.Range('E30').Value = '  +0.15%  '
$ws.Range('B31').Value = 'Hedera'
$ws.Range('C31').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0552'
$ws.Range('E31').Value = '  +0.89%  '
$ws.Range('B32').Value = 'Filecoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.99'
$ws.Range('E32').Value = '  +2.92%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.96'
$ws.Range('E33').Value = '  -0.45%  '
$ws.Range('E34').Value = '  +11.22%  '
$ws.Range('E35').Value = '  +4.44%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.700'
$ws.Range('E36').Value = '  +3.71%  '
$ws.Range('B37').Value = 'WEMIXToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.39'
$ws.Range('E37').Value = '  +9.14%  '
$ws.Range('B38').Value = 'Aave'
$ws.Range('C38').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '93.21'
$ws.Range('E38').Value = '  +0.04%  '
$ws.Range('D39').Value = '1.337.91'
$ws.Range('E39').Value = '  +2.30%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0195'
$ws.Range('E40').Value = '  +1.71%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.995'
$ws.Range('E41').Value = '  +1.50%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '14.82'
$ws.Range('E42').Value = '  +0.98%  '
$ws.Range('E43').Value = '  -1.30%  '
$ws.Range('E44').Value = '  -0.80%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.77'
$ws.Range('E45').Value = '  +0.09%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '6.29'
$ws.Range('E46').Value = '  +1.33%  '
$ws.Range('E47').Value = '  +2.08%  '
$ws.Range('D48').Value = '2.008.27'
$ws.Range('E48').Value = '  +0.80%  '
$ws.Range('E49').Value = '  +0.05%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0669'
$ws.Range('E50').Value = '  +4.89%  '
$ws.Range('E51').Value = '  +13.86%  '
